# Timesheet update: add Feb 19 2020 entries (rows 162-171) to Sheet1.
# Mirrors the author's entry order so the shared-strings table is
# rebuilt in the same sequence as the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 162: blank separator row (same highlighted style as row 151) ---
$sepA = $ws.Cells.Item(162, 1)
$sepA.Interior.Color = 49407
$sepA.HorizontalAlignment = -4108
$sepA.VerticalAlignment = -4108

$sepB = $ws.Cells.Item(162, 2)
$sepB.Interior.Color = 49407
$sepB.HorizontalAlignment = -4131
$sepB.VerticalAlignment = -4108

$sepC = $ws.Cells.Item(162, 3)
$sepC.Interior.Color = 49407
$sepC.HorizontalAlignment = -4108
$sepC.VerticalAlignment = -4108

# --- Row 163 ---
$ws.Cells.Item(163, 1).Value = "Feb 19 10:00 to 11:00"
$b163 = $ws.Cells.Item(163, 2)
$b163.Value = "Saved transformed data, used json configuration file, but logic did not worked. Saved`r`nusing declaring path in class."
$b163.WrapText = $true
$b163.HorizontalAlignment = -4131
$b163.VerticalAlignment = -4108
$ws.Rows.Item(163).RowHeight = 30
$ws.Cells.Item(163, 3).Value = "Infimetrics"

# --- Row 164 ---
$ws.Cells.Item(164, 1).Value = "Feb 19 11:00 to 12:00"
$ws.Cells.Item(164, 2).Value = "Documenting code properly"
$ws.Cells.Item(164, 3).Value = "Infimetrics"

# --- Row 165 ---
$ws.Cells.Item(165, 1).Value = "Feb 19 12:00 to 13:00"
$ws.Cells.Item(165, 3).Value = "Infimetrics"

# --- Row 166 ---
$ws.Cells.Item(166, 1).Value = "Feb 19 13:00 to 14:00"
$ws.Cells.Item(166, 3).Value = "Infimetrics"

# Descriptions for rows 165/166 filled in after both timestamps (165, 166)
# were entered - matches the shared-string order in the source file.
$ws.Cells.Item(165, 2).Value = "Documented and commented data transformation"
$ws.Cells.Item(166, 2).Value = "Worked on code generator problem"

# --- Row 167 ---
$ws.Cells.Item(167, 1).Value = "Feb 19 14:00 to 15:00"
$ws.Cells.Item(167, 2).Value = "Lunch"
$ws.Cells.Item(167, 3).Value = "Infimetrics"

# --- Row 168 ---
$ws.Cells.Item(168, 1).Value = "Feb 19 15:00 to 16:00"
$ws.Cells.Item(168, 2).Value = "Started statistical analysis phase of ml pipeline. Created new jupyter notebook"
$ws.Cells.Item(168, 3).Value = "Infimetrics"

# --- Row 169 ---
$ws.Cells.Item(169, 1).Value = "Feb 19 16:00 to 17:00"
$ws.Cells.Item(169, 2).Value = "Checking null values"
$ws.Cells.Item(169, 3).Value = "Infimetrics"

# --- Row 170 ---
$ws.Cells.Item(170, 1).Value = "Feb 19 17:00 to 18:00"
$ws.Cells.Item(170, 2).Value = "Printing total null values of both data. Displayed descriptive statistic of selected columns"
$ws.Cells.Item(170, 3).Value = "Infimetrics"

# --- Row 171 ---
$ws.Cells.Item(171, 1).Value = "Feb 19 18:00 to 19:00"
$b171 = $ws.Cells.Item(171, 2)
$b171.Value = "Identified outliers in alarm and automation duration. Working on outlier removal, lag`r`ndetections and statistical hypothesis testing."
$b171.WrapText = $true
$b171.HorizontalAlignment = -4131
$b171.VerticalAlignment = -4108
$ws.Rows.Item(171).RowHeight = 30
$ws.Cells.Item(171, 3).Value = "Infimetrics"

# --- View state: selection + active cell move to the new last row ---
$ws.Range("D171").Select()
$excel.ActiveWindow.ScrollRow = 157
